# Aula 001 - slide edits
# 1) presentation.xml gains an empty p15:sldGuideLst extension (added by
#    PowerPoint automatically once the guides UI is touched / file resaved).
# 2) Slide 1: resize + retitle the subtitle placeholder.
# 3) Slide 5: shrink the content placeholder and add a new "GitHub" textbox
#    below it (duplicated from the placeholder so it inherits its style).

# Helper: PowerPoint's Left/Top/Width/Height setters take points and the
# points -> EMU conversion used by this host truncates, so nudge by half an
# EMU (in point units) before assigning to land on the exact target EMU.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation

# --- Slide 1: subtitle shape -------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(1)
$subtitle.Height = EmuToPt 503618
$subtitle.TextFrame.TextRange.Text = "Leonardo Barcelos Marques"

# --- Slide 5: content placeholder + new shape --------------------------------
$s5 = $p.Slides.Item(5)
$placeholder = $s5.Shapes.Item(3)
$placeholder.Height = EmuToPt 1685671

# Duplicate the placeholder so the new shape inherits its list style / body
# formatting, then reposition it below and replace its text.
$dupRange = $placeholder.Duplicate()
$newShape = $dupRange.Item(1)
$newShape.Left = EmuToPt 838200
$newShape.Top = EmuToPt 3523361
$newShape.Width = EmuToPt 10515600
$newShape.Height = EmuToPt 1685671

$tr = $newShape.TextFrame.TextRange
$tr.Text = "GitHub`rhttps://github.com/LBarcelosM/Edge-Blazor-Course.git"
